# User-Journey-Map-Portfolio-Project.pptx
# "Add peer review statements & update files"
#
# The journey-map table on slide 1 was reflowed: the table grew a bit
# taller/narrower, moved up slightly, and the "TASK LIST" row (row 2)
# grew to make room for the added text. Reproduce that relayout via the
# Table/Shape COM object model (EMU -> points, 1 pt = 12700 EMU).

$EMU_PER_PT = 12700

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the journey-map table on the slide (rather than hard-coding an index).
$tableShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasTable) {
        $tableShape = $candidate
    }
}

$tbl = $tableShape.Table

# New column widths (EMU), left to right.
$colWidths = @(1663200, 1663200, 1663200, 1639440, 1685160, 1664280)
for ($c = 1; $c -le $colWidths.Count; $c++) {
    $tbl.Columns.Item($c).Width = $colWidths[$c - 1] / $EMU_PER_PT
}

# Row 2 ("TASK LIST") grows from 1850040 EMU to 2013120 EMU; the other
# rows keep their original heights.
$tbl.Rows.Item(2).Height = 2013120 / $EMU_PER_PT

# The table frame shifts up (its x position and its width/height are
# driven by the column/row edits above); nudge Top to the new value.
# (A hair off the plain EMU/12700 quotient to land exactly on 1341360 EMU
# once PowerPoint's internal point-precision rounding is applied.)
$tableShape.Top = 105.61894763779527
